$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "col2" header between the existing col1 and col3 headers,
# and add two new data rows -- turning the original 1-row dataframe
# (col1,col3 / a3,4) into a 3-row x 3-col dataframe with an extra
# middle column (col2) filled in for every row.

# Header row
$ws.Range("A1").Value = "col1"
$ws.Range("B1").Value = "col2"
$ws.Range("C1").Value = "col3"

# Data rows
$ws.Range("A2").Value = "a1"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3.1

$ws.Range("A3").Value = "a3"
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 5.1

# A4 ("6") must land in the sheet as literal text, not a number -- a plain
# Range.Value assignment of a numeric-looking string auto-coerces to a
# number (and forcing text via NumberFormat="@" or an apostrophe prefix
# would stamp an extra cell style that the original file doesn't have).
# Route it through a text formula + paste-values-only instead, which
# keeps the workbook's cell styles untouched.
$ws.Range("Z1").Formula = "=""6"""
$ws.Range("Z1").Copy()
$ws.Range("A4").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

$ws.Range("B4").Value = 7
$ws.Range("C4").Value = 8.1
